# Update NATMI ligand/receptor expression & specificity metrics with
# recomputed TPM-based values (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("O2").Value = 0.5024565340298992
$ws.Range("P2").Value = 0.602355957672625
$ws.Range("Q2").Value = 16.59853721842
$ws.Range("R2").Value = 149.38683496578
$ws.Range("S2").Value = 0.122099265911213
$ws.Range("T2").Value = 0.1511989845835743
$ws.Range("G3").Value = 144.783305
$ws.Range("H3").Value = 434.349915
$ws.Range("I3").Value = 0.2430046335191003
$ws.Range("J3").Value = 0.251012682214973
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.113523
$ws.Range("N3").Value = 0.227046
$ws.Range("O3").Value = 0.4975434659701009
$ws.Range("P3").Value = 0.3976440423273752
$ws.Range("Q3").Value = 16.436235133515
$ws.Range("R3").Value = 98.61741080108999
$ws.Range("S3").Value = 0.1209053676078873
$ws.Range("T3").Value = 0.09981369763139868
$ws.Range("G4").Value = 82.24887099999999
$ws.Range("I4").Value = 0.1380466950572427
$ws.Range("J4").Value = 0.1425959278859072
$ws.Range("O4").Value = 0.5024565340298992
$ws.Range("P4").Value = 0.602355957672625
$ws.Range("Q4").Value = 9.429339566924
$ws.Range("R4").Value = 84.864056102316
$ws.Range("S4").Value = 0.06936246393274458
$ws.Range("T4").Value = 0.08589350670193219
$ws.Range("G5").Value = 82.24887099999999
$ws.Range("I5").Value = 0.1380466950572427
$ws.Range("J5").Value = 0.1425959278859072
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.113523
$ws.Range("N5").Value = 0.227046
$ws.Range("O5").Value = 0.4975434659701009
$ws.Range("P5").Value = 0.3976440423273752
$ws.Range("Q5").Value = 9.337138582532999
$ws.Range("R5").Value = 56.022831495198
$ws.Range("S5").Value = 0.06868423112449812
$ws.Range("T5").Value = 0.05670242118397501
$ws.Range("G6").Value = 163.8590903333333
$ws.Range("H6").Value = 491.577271
$ws.Range("I6").Value = 0.2750214756820535
$ws.Range("J6").Value = 0.284084617144743
$ws.Range("O6").Value = 0.5024565340298992
$ws.Range("P6").Value = 0.602355957672625
$ws.Range("Q6").Value = 18.78546155217467
$ws.Range("R6").Value = 169.069153969572
$ws.Range("S6").Value = 0.1381863374549928
$ws.Range("T6").Value = 0.1711200616202827
$ws.Range("G7").Value = 163.8590903333333
$ws.Range("H7").Value = 491.577271
$ws.Range("I7").Value = 0.2750214756820535
$ws.Range("J7").Value = 0.284084617144743
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.113523
$ws.Range("N7").Value = 0.227046
$ws.Range("O7").Value = 0.4975434659701009
$ws.Range("P7").Value = 0.3976440423273752
$ws.Range("Q7").Value = 18.601775511911
$ws.Range("R7").Value = 111.610653071466
$ws.Range("S7").Value = 0.1368351382270607
$ws.Range("T7").Value = 0.1129645555244603
$ws.Range("G8").Value = 57.0238095
$ws.Range("H8").Value = 114.047619
$ws.Range("I8").Value = 0.09570889357312636
$ws.Range("J8").Value = 0.06590860906562239
$ws.Range("O8").Value = 0.5024565340298992
$ws.Range("P8").Value = 0.602355957672625
$ws.Range("Q8").Value = 6.537437616318001
$ws.Range("R8").Value = 39.224625697908
$ws.Range("S8").Value = 0.04808955894058956
$ws.Range("T8").Value = 0.03970044333259363
$ws.Range("G9").Value = 57.0238095
$ws.Range("H9").Value = 114.047619
$ws.Range("I9").Value = 0.09570889357312636
$ws.Range("J9").Value = 0.06590860906562239
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.5
$ws.Range("M9").Value = 0.113523
$ws.Range("N9").Value = 0.227046
$ws.Range("O9").Value = 0.4975434659701009
$ws.Range("P9").Value = 0.3976440423273752
$ws.Range("Q9").Value = 6.4735139258685
$ws.Range("R9").Value = 25.894055703474
$ws.Range("S9").Value = 0.0476193346325368
$ws.Range("T9").Value = 0.02620816573302877
$ws.Range("G10").Value = 147.8896333333333
$ws.Range("H10").Value = 443.6689
$ws.Range("I10").Value = 0.2482183021684772
$ws.Range("J10").Value = 0.2563981636887546
$ws.Range("O10").Value = 0.5024565340298992
$ws.Range("P10").Value = 0.602355957672625
$ws.Range("Q10").Value = 16.95465912386667
$ws.Range("R10").Value = 152.5919321148
$ws.Range("S10").Value = 0.1247189077903593
$ws.Range("T10").Value = 0.1544429614342422
$ws.Range("G11").Value = 147.8896333333333
$ws.Range("H11").Value = 443.6689
$ws.Range("I11").Value = 0.2482183021684772
$ws.Range("J11").Value = 0.2563981636887546
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.5
$ws.Range("M11").Value = 0.113523
$ws.Range("N11").Value = 0.227046
$ws.Range("O11").Value = 0.4975434659701009
$ws.Range("P11").Value = 0.3976440423273752
$ws.Range("Q11").Value = 16.7888748449
$ws.Range("R11").Value = 100.7332490694
$ws.Range("S11").Value = 0.123499394378118
$ws.Range("T11").Value = 0.1019552022545124
